# "week 3 data updated" -- fill in previously-missing cortisol collection
# times for a few participant/week rows, fix up the "participant forgot..."
# note text, and update the view/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$noteCogTask = "participant forgot to do the cortisol collection the day after the cognitive task, so she collected it 2 days after the cog task instead (final time is technically the next day)"
$noteNextDay = "(final time is technically the next day)"

# --- Row 16: Starlight, Week 3 ---
$ws.Cells.Item(16, 3).Value = 0.44027777777777777
$ws.Cells.Item(16, 4).Value = 0.47152777777777777
$ws.Cells.Item(16, 5).Value = 0.69027777777777777
$ws.Cells.Item(16, 6).Value = 0.94027777777777777
$ws.Cells.Item(16, 7).Value = 0.05486111111111111
$ws.Range("C16:G16").NumberFormat = "h:mm"
$ws.Cells.Item(16, 8).Value = $noteCogTask

# --- Row 24: Prism, Week 2 ---
$ws.Cells.Item(24, 3).Value = 0.33333333333333331
$ws.Cells.Item(24, 4).Value = 0.36527777777777776
$ws.Cells.Item(24, 5).Value = 0.67152777777777772
$ws.Cells.Item(24, 6).Value = 0.85416666666666663
$ws.Cells.Item(24, 7).Value = 0.97916666666666663
$ws.Range("C24:G24").NumberFormat = "h:mm"

# --- Row 27: Cascade, Week 2 -- correct the wording of the existing note ---
$ws.Cells.Item(27, 8).Value = $noteCogTask

# --- Row 37: Quartz, Week 3 ---
$ws.Cells.Item(37, 3).Value = 0.32291666666666669
$ws.Cells.Item(37, 4).Value = 0.35416666666666669
$ws.Cells.Item(37, 5).Value = 0.65625
$ws.Cells.Item(37, 6).Value = 0.82291666666666663
$ws.Cells.Item(37, 7).Value = 0.07291666666666667
$ws.Range("C37:G37").NumberFormat = "h:mm"
$ws.Cells.Item(37, 8).Value = $noteNextDay

# --- Column H grew wider to fit the longer note text ---
$ws.Columns.Item(8).ColumnWidth = 145.6328125

# --- Update the active selection/view state ---
$ws.Range("H9").Select()
